$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A95").Value = 20220831
$ws.Range("B95").Value = 2210.5472799999998
$ws.Range("C95").Value = 2224.4699999999998
$ws.Range("D95").Formula = "=100*(B95-C95)/C95"
$ws.Range("E95").Value = 180
$ws.Range("F95").Value = "CRM OPENED 20220825 LHZ"

$ws.Range("G96").Select()
